$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '329.86'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.18%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '3.67%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.650'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-4.00%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08161'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.72%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.036'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '5.18%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.749'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.13%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.949'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.09%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9185'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.45%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1258'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.72%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1953'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.06%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09366'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2.37%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03693'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '4.58%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1055'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '10.28%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001301'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.64%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006171'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.10%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.430'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '2.33%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.544'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.59%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.276'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-5.25%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.71%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '10.21%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04409'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.03%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001270'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.69%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004301'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-2.02%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '3.78%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02760'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '13.90%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05446'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '4.11%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007657'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.74%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009478'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.13%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1416'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.70%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002115'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.21%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01198'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '6.84%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006881'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '2.38%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.25%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '60.54%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.003581'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '19.34%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.25%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002004'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.25%'
